# "Update countries & provincias Spain" - refresh the COVID-19 country
# data snapshot (header timestamp + the rows whose counters moved since
# the last pull). Countries keep their original row only when their
# "Casos totales" rank doesn't change; Colombia/Peru, Haiti/Tunez,
# Luxemburgo/Zimbabue and Montserrat/Islas Malvinas swap rows because
# their updated totals re-order them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{}
$rows[1] = @('Datos actualizados a 19 de Septiembre de 2020 a las 00:20', $null, $null, $null, $null, $null, $null, $null)
$rows[4] = @('Estados Unidos', 6918358, 43762, 4182446, 2532846, 0, 853, 203066)
$rows[6] = @('Brasil', 4495183, 37740, 3753082, 606308, 0, 762, 135793)
$rows[8] = @('Colombia', 750471, 6526, 621521, 105100, 0, 185, 23850)
$rows[9] = @('Peru', 750098, 0, 594513, 124439, 0, 0, 31146)
$rows[25] = @('Alemania', 271244, 2202, 241300, 20480, 0, 7, 9464)
$rows[45] = @('Guatemala', 84344, 680, 73748, 7520, 0, 40, 3076)
$rows[61] = @('Suiza', 49283, 488, 40500, 6738, 0, 3, 2045)
$rows[63] = @('Ghana', 45760, 46, 44973, 492, 0, 1, 295)
$rows[82] = @('Camerun', 20371, 68, 19124, 831, 0, 1, 416)
$rows[84] = @('Bulgaria', 18733, 189, 13510, 4470, 0, 4, 753)
$rows[91] = @('Zambia', 14022, 94, 13207, 486, 0, 3, 329)
$rows[104] = @('Haiti', 8600, 44, 6363, 2016, 0, 1, 221)
$rows[105] = @('Tunez', 8570, 0, 2342, 6095, 0, 0, 133)
$rows[106] = @('Luxemburgo', 7718, 177, 6703, 891, 0, 0, 124)
$rows[107] = @('Zimbabue', 7647, 14, 5883, 1540, 0, 0, 224)
$rows[115] = @('Suazilandia', 5215, 24, 4478, 634, 0, 0, 103)
$rows[136] = @('Aruba', 3460, 78, 2128, 1309, 0, 0, 23)
$rows[137] = @('Guadalupe', 3426, 0, 837, 2563, 0, 0, 26)
$rows[138] = @('Somalia', 3390, 0, 2812, 480, 0, 0, 98)
$rows[146] = @('Sudan del Sur', 2609, 10, 1290, 1270, 0, 0, 49)
$rows[152] = @('Guyana', 2102, 75, 1314, 726, 0, 2, 62)
$rows[156] = @('Burkina Faso', 1797, 30, 1173, 568, 0, 0, 56)
$rows[157] = @('Togo', 1640, 22, 1251, 348, 0, 0, 41)
$rows[169] = @('Santo Tome y Principe', 908, 1, 873, 20, 0, 0, 15)
$rows[191] = @('Barbados', 185, 0, 172, 6, 0, 0, 7)
$rows[197] = @('Islas Virgenes Britanicas', 69, 0, 48, 20, 0, 0, 1)
$rows[214] = @('Montserrat', 13, 0, 12, 0, 0, 0, 1)
$rows[215] = @('Islas Malvinas', 13, 0, 13, 0, 0, 0, 0)

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $c = 1
    foreach ($v in $vals) {
        $ws.Cells.Item($r, $c).Value = $v
        $c = $c + 1
    }
}
